$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : NumPos goes from 2 -> 1 ("cheminData" parameter, values unchanged) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "cheminData"
$ws.Range("C2").Value = "/data"
$ws.Range("D2").Value = "path fichiers data"
$ws.Range("E2").Value = 44257

# --- Row 3 : NumPos goes from 3 -> 2, "cheminUsers" renamed "cheminUtilisateurs" ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "cheminUtilisateurs"
$ws.Range("C3").Value = "/utilisateurs"
$ws.Range("D3").Value = "path fichiers users"
$ws.Range("E3").Value = 44250

# --- Row 4 (new) : "cheminMessages" parameter, inserted before the old row 4 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "cheminMessages"
$ws.Range("C4").Value = "/data"
$ws.Range("D4").Value = "path fichier messages.xlsx"
$ws.Range("E4").Value = 44258

# --- Row 5 : former row 4 ("imprimeOK"), now shifted down one row ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "imprimeOK"
$ws.Range("C5").Value = $true
$ws.Range("D5").Value = "Attention ! Valeurs booléennes (en français ici)."
$ws.Range("E5").Value = 44254
# Re-apply the short-date format (the cell previously held the blank
# placeholder row, which used the long date-time format).
$ws.Range("E2").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 8 : E8 switches to the datetime number format ---
$ws.Range("E8").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"

# --- Row 20 (new trailing blank row, mirrors row 19's formatting) ---
$ws.Range("A19:E19").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Range("A20:E20").ClearContents()
$excel.CutCopyMode = 0

# --- Column width adjustments (columns B and C grew to fit the new text) ---
$ws.Columns("B").ColumnWidth = 15.3
$ws.Columns("C").ColumnWidth = 39

# --- Selection as left by the author ---
$ws.Range("B4").Select()
